$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "37.098.59"
$ws.Range('E2').Value = "  +1.49%  "
$ws.Range('D3').Value = "2.059.57"
$ws.Range('E3').Value = "  -2.14%  "
$ws.Range('E4').Value = "  +0.04%  "
$ws.Range('D5').Value = "'249.01"
$ws.Range('E5').Value = "  -1.22%  "
$ws.Range('D6').Value = "'0.659"
$ws.Range('E6').Value = "  -0.12%  "
$ws.Range('E7').Value = "  -0.08%  "
$ws.Range('D8').Value = "'55.48"
$ws.Range('E8').Value = "  +16.95%  "
$ws.Range('D9').Value = "'60.97"
$ws.Range('E9').Value = "  +2.87%  "
$ws.Range('E10').Value = "  +1.84%  "
$ws.Range('E11').Value = "  +7.34%  "
$ws.Range('E12').Value = "  +5.82%  "
$ws.Range('D13').Value = "'15.10"
$ws.Range('E13').Value = "  +5.91%  "
$ws.Range('D14').Value = "2.358.41"
$ws.Range('E14').Value = "  -2.32%  "
$ws.Range('E15').Value = "  -1.17%  "
$ws.Range('D16').Value = "'5.25"
$ws.Range('E16').Value = "  +3.31%  "
$ws.Range('D17').Value = "2.065.92"
$ws.Range('E17').Value = "  -1.98%  "
$ws.Range('D18').Value = "37.041.72"
$ws.Range('E18').Value = "  +1.31%  "
$ws.Range('D19').Value = "0.0₃0957"
$ws.Range('E19').Value = "  +15.47%  "
$ws.Range('D20').Value = "'72.51"
$ws.Range('E20').Value = "  -0.49%  "
$ws.Range('D21').Value = "'14.21"
$ws.Range('E21').Value = "  +7.52%  "
$ws.Range('E22').Value = "  +4.27%  "
$ws.Range('D23').Value = "'237.33"
$ws.Range('E23').Value = "  -0.98%  "
$ws.Range('E24').Value = "  +0.01%  "
$ws.Range('E25').Value = "  -0.78%  "
$ws.Range('D26').Value = "'170.66"
$ws.Range('E26').Value = "  -0.44%  "
$ws.Range('D27').Value = "'9.06"
$ws.Range('E27').Value = "  -0.85%  "
$ws.Range('D28').Value = "'20.13"
$ws.Range('E28').Value = "  -5.63%  "
$ws.Range('D29').Value = "'1.98"
$ws.Range('E29').Value = "  +0.31%  "
$ws.Range('E30').Value = "  +0.34%  "
$ws.Range('E31').Value = "  +2.85%  "
$ws.Range('E32').Value = "  +11.23%  "
$ws.Range('D33').Value = "'0.0624"
$ws.Range('E33').Value = "  +3.03%  "
$ws.Range('E34').Value = "  +7.53%  "
$ws.Range('E35').Value = "  -0.05%  "
$ws.Range('D36').Value = "'2.30"
$ws.Range('E36').Value = "  -0.74%  "
$ws.Range('D37').Value = "'0.0855"
$ws.Range('E37').Value = "  -3.47%  "
$ws.Range('D38').Value = "'1.77"
$ws.Range('E38').Value = "  -6.16%  "
$ws.Range('E39').Value = "  +1.51%  "
$ws.Range('D40').Value = "'0.106"
$ws.Range('E40').Value = "  +26.68%  "
$ws.Range('D41').Value = "'18.10"
$ws.Range('E41').Value = "  +11.59%  "
$ws.Range('E42').Value = "  +0.37%  "
$ws.Range('E43').Value = "  -2.75%  "
$ws.Range('D44').Value = "'96.54"
$ws.Range('E44').Value = "  -1.14%  "
$ws.Range('B45').Value = "HuobiToken"
$ws.Range('C45').Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range('D45').Value = "'2.77"
$ws.Range('E45').Value = "  +0.72%  "
$ws.Range('B46').Value = "FTXToken"
$ws.Range('C46').Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range('D46').Value = "'4.21"
$ws.Range('E46').Value = "  +47.77%  "
$ws.Range('E47').Value = "  +8.07%  "
$ws.Range('D48').Value = "'13.23"
$ws.Range('E48').Value = "  -52.82%  "
$ws.Range('D49').Value = "1.297.73"
$ws.Range('E49').Value = "  -2.68%  "
$ws.Range('D50').Value = "'2.92"
$ws.Range('E50').Value = "  +2.82%  "
$ws.Range('B51').Value = "THORChain"
$ws.Range('C51').Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range('D51').Value = "'4.05"
$ws.Range('E51').Value = "  +6.09%  "
